$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.071.79'
$ws.Range("E2").Value = '  -2.12%  '
$ws.Range("D3").Value = '1.826.43'
$ws.Range("E3").Value = '  -0.99%  '
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = '  -0.83%  '
$ws.Range("D5").Value = "'311.58"
$ws.Range("E5").Value = '  -1.85%  '
$ws.Range("E6").Value = '  -0.70%  '
$ws.Range("D7").Value = "'0.4240"
$ws.Range("E7").Value = '  -1.27%  '
$ws.Range("D8").Value = "'0.3672"
$ws.Range("E8").Value = '  -1.79%  '
$ws.Range("D9").Value = "'0.07231"
$ws.Range("E9").Value = '  -1.16%  '
$ws.Range("D10").Value = "'0.8462"
$ws.Range("E10").Value = '  -3.05%  '
$ws.Range("D11").Value = "'20.70"
$ws.Range("E11").Value = '  -3.46%  '
$ws.Range("D12").Value = '1.825.77'
$ws.Range("E12").Value = '  -1.05%  '
$ws.Range("D13").Value = "'6.665"
$ws.Range("E13").Value = '  -0.76%  '
$ws.Range("E14").Value = '  -2.13%  '
$ws.Range("D15").Value = "'0.07040"
$ws.Range("E15").Value = '  -0.96%  '
$ws.Range("D16").Value = "'89.66"
$ws.Range("E16").Value = '  +1.02%  '
$ws.Range("E17").Value = '  -0.82%  '
$ws.Range("D18").Value = "'0.000008763"
$ws.Range("E18").Value = '  -2.50%  '
$ws.Range("E19").Value = '  -0.64%  '
$ws.Range("E20").Value = '  -3.31%  '
$ws.Range("D21").Value = '27.119.80'
$ws.Range("E21").Value = '  -1.96%  '
$ws.Range("D22").Value = "'5.135"
$ws.Range("E22").Value = '  -1.16%  '
$ws.Range("E23").Value = '  -1.77%  '
$ws.Range("D24").Value = '2.050.87'
$ws.Range("D25").Value = "'1.977"
$ws.Range("E25").Value = '  +0.34%  '
$ws.Range("D26").Value = "'151.53"
$ws.Range("E26").Value = '  -2.09%  '
$ws.Range("D27").Value = "'2.251"
$ws.Range("E27").Value = '  +4.42%  '
$ws.Range("E28").Value = '  -1.65%  '
$ws.Range("D29").Value = "'5.254"
$ws.Range("E29").Value = '  -1.43%  '
$ws.Range("D30").Value = "'116.84"
$ws.Range("E30").Value = '  -0.77%  '
$ws.Range("D31").Value = "'0.08715"
$ws.Range("E31").Value = '  -2.14%  '
$ws.Range("D32").Value = "'1.181"
$ws.Range("E32").Value = '  -3.22%  '
$ws.Range("D33").Value = "'0.7376"
$ws.Range("E33").Value = '  -4.71%  '
$ws.Range("E34").Value = '  +0.22%  '
$ws.Range("D35").Value = "'4.431"
$ws.Range("E35").Value = '  -2.10%  '
$ws.Range("E36").Value = '  -0.82%  '
$ws.Range("D37").Value = "'1.093"
$ws.Range("D38").Value = "'0.01945"
$ws.Range("E38").Value = '  -1.70%  '
$ws.Range("D39").Value = "'0.05242"
$ws.Range("E39").Value = '  -1.21%  '
$ws.Range("D40").Value = "'7.338"
$ws.Range("E40").Value = '  +2.74%  '
$ws.Range("D41").Value = "'2.873"
$ws.Range("E41").Value = '  -0.40%  '
$ws.Range("D42").Value = "'0.1689"
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").Value = "'0.5073"
$ws.Range("E43").Value = '  -0.99%  '
$ws.Range("D44").Value = "'8.567"
$ws.Range("E44").Value = '  -2.28%  '
$ws.Range("D45").Value = "'1.986"
$ws.Range("E45").Value = '  +7.89%  '
$ws.Range("D46").Value = "'10.49"
$ws.Range("E46").Value = '  -1.41%  '
$ws.Range("D47").Value = "'0.4740"
$ws.Range("E47").Value = '  -0.10%  '
$ws.Range("D48").Value = "'105.77"
$ws.Range("E48").Value = '  -1.40%  '
$ws.Range("D49").Value = "'1.001"
$ws.Range("E49").Value = '  -0.81%  '
$ws.Range("D50").Value = "'0.06324"
$ws.Range("E50").Value = '  -1.91%  '
$ws.Range("D51").Value = "'1.651"

# Strip the auto-applied "Text" number format (quote-prefix) so cell
# styling matches the source file, which had no explicit style override.
foreach ($addr in @("D4","D5","D7","D8","D9","D10","D11","D13","D15","D16","D18","D22","D25","D26","D27","D29","D30","D31","D32","D33","D35","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")) {
    $ws.Range($addr).ClearFormats()
}
